# Weekly update: a new price-report row is inserted at row 3 (just below the
# header/first data row) and every existing data row shifts down by one.
# The last existing row (old row 32) ends up as the new row 33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 3 - this pushes rows 3..32 down to 4..33
# and Excel naturally carries the existing row formatting (incl. the date
# style on column D) down onto the new row, matching the target OOXML.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with this week's record. Columns that are
# constant across every row in this table (A, B, C, E, F, G, H, I, N, O, Q, R)
# keep the same values as the rest of the sheet.
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 44819
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 100112043
$ws.Range("G3").Value = "Pepino dulce"
$ws.Range("H3").Value = "Cultivar IV Región"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 15000
$ws.Range("N3").Value = "$/bandeja 18 kilos"
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 833
$ws.Range("Q3").Value = 18
$ws.Range("R3").Value = "Hortaliza"
